# Initial Data File Updated, Debt table and calculations added

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Rename the original sheet
$ws1.Name = "Transacciones"

# 2) Add the new "Deudas" sheet right after "Transacciones"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Deudas"

# 3) Shift the transaction dates in "Transacciones" forward by one month (Mar -> Apr)
$ws1.Range("A7").Value = 43556
$ws1.Range("A8").Value = 43556
$ws1.Range("A9").Value = 43557
$ws1.Range("A10").Value = 43557
$ws1.Range("A11").Value = 43557
$ws1.Range("A12").Value = 43558
$ws1.Range("A13").Value = 43558
$ws1.Range("A14").Value = 43558

# 4) Add the new row 15 transaction (electricity bill payment)
$ws1.Range("A15").Value = 43559
$ws1.Range("B15").Value = 63
$ws1.Range("C15").Value = "Pago de Recibo de Electricidad"
$ws1.Range("D15").Value = "Servicios"
$ws1.Range("E15").Value = "Gasto"
$ws1.Range("F15").Value = "Tarjeta Santander"
$ws1.Range("G15").Value = "Transferencia"
$ws1.Range("K15").Value = 7831.82
$ws1.Range("L15").Formula = "=L14-B15"
$ws1.Range("M15").Value = 70
$ws1.Range("N15").Formula = "=SUM(K15:M15)"
$ws1.Range("O15").Formula = "=N15-4000"

# 5) Update the running-balance block (S/T/U/V columns)
$ws1.Range("S15").Value = 3988.82
$ws1.Range("T15").Value = 250
$ws1.Range("U15").Formula = "=S15-T19"
$ws1.Range("V15").Formula = "=U15+U12"
$ws1.Range("T17").Value = 0

# 6) Populate the new "Deudas" sheet
$ws2.Range("F2").Value = "Deuda Rentas"
$ws2.Range("G2").Formula = "=SUM(`$B4:`$B8)"
$ws2.Range("M2").Value = "Deuda Restante (Rentas)"
$ws2.Range("N2").Formula = "=IF((G2-SUM(J:J)) <= 0,0,G2-SUM(J:J))"
$ws2.Range("P2").Value = "Pago mensual"
$ws2.Range("Q2").Value = 3000

$ws2.Range("A3").Value = "Fecha"
$ws2.Range("B3").Value = "Monto"
$ws2.Range("C3").Value = "Descripción"
$ws2.Range("F3").Value = "Deuda Total"
$ws2.Range("G3").Formula = "=SUM(B:B)"
$ws2.Range("I3").Value = "Fecha"
$ws2.Range("J3").Value = "Monto "
$ws2.Range("K3").Value = "Descripción"
$ws2.Range("M3").Value = "Deuda Restante (Total)"
$ws2.Range("N3").Formula = "=G3-SUM(J:J)"
$ws2.Range("P3").Value = "Estimado Meses restantes"
$ws2.Range("Q3").Formula = "=N3/Q2"

$ws2.Range("B4").Value = 2500
$ws2.Range("C4").Value = "Renta 1er departamento"
$ws2.Range("I4").Value = 43422
$ws2.Range("J4").Value = 2400
$ws2.Range("K4").Value = "Pago"

$ws2.Range("B5").Value = 2500
$ws2.Range("C5").Value = "Depósito 1er departamento"
$ws2.Range("I5").Value = 43449
$ws2.Range("J5").Value = 2000

$ws2.Range("B6").Value = 4000
$ws2.Range("C6").Value = "1er Mes de Manutención"
$ws2.Range("I6").Value = 43462
$ws2.Range("J6").Value = 1500

$ws2.Range("B7").Value = 4900
$ws2.Range("C7").Value = "Renta 2do Departamento"
$ws2.Range("I7").Value = 43480
$ws2.Range("J7").Value = 1500

$ws2.Range("B8").Value = 2000
$ws2.Range("C8").Value = "Depósito 2do Departamento"
$ws2.Range("I8").Value = 43496
$ws2.Range("J8").Value = 1500

$ws2.Range("B9").Value = 11759
$ws2.Range("C9").Value = "Xbox One X"
$ws2.Range("I9").Value = 43511
$ws2.Range("J9").Value = 1500

$ws2.Range("B10").Value = 20324
$ws2.Range("C10").Value = "iPhone 8 Plus"
$ws2.Range("I10").Value = 43524
$ws2.Range("J10").Value = 1500

$ws2.Range("B11").Value = 10500
$ws2.Range("C11").Value = "TV Sharp 4K 50`""
$ws2.Range("I11").Value = 43539
$ws2.Range("J11").Value = 1500

$ws2.Range("I12").Value = 43553
$ws2.Range("J12").Value = 1500

# Apply the date number format (matches column A / I date cells elsewhere) to the date columns
$ws1.Range("A15").NumberFormat = $ws1.Range("A14").NumberFormat
$ws2.Range("I4:I12").NumberFormat = $ws1.Range("A14").NumberFormat

# 7) Restore the original active sheet / selection state
$ws1.Activate()
$ws1.Range("S16").Select()
